$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 -> add_user ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "add_user"

# --- Pre-format the new ranges so new cells pick up the existing style (s="1")
#     without creating a brand-new cellXf entry. Only format the cells that
#     will actually receive a value, so cells that stay blank in the target
#     (C3, D4, B6 on the add_user sheet) are never materialised. ---
$ws1.Range("A1").Copy()
$ws1.Range("D1:D8").PasteSpecial(-4122)

$ws1.Range("A1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)
$ws2.Range("A2:G2").PasteSpecial(-4122)
$ws2.Range("A3:B3").PasteSpecial(-4122)
$ws2.Range("D3:G3").PasteSpecial(-4122)
$ws2.Range("A4:C4").PasteSpecial(-4122)
$ws2.Range("E4:G4").PasteSpecial(-4122)
$ws2.Range("A5:G5").PasteSpecial(-4122)
$ws2.Range("A6").PasteSpecial(-4122)
$ws2.Range("C6:G6").PasteSpecial(-4122)
$ws2.Range("A7:G7").PasteSpecial(-4122)

# ============ Sheet1 (login) : new column D "testcase name" ============
$ws1.Range("D1").Value = "testcase name"
$ws1.Range("D2").Value = "Valid testcase"
$ws1.Range("D3").Value = "Wrong password"
$ws1.Range("D4").Value = "Non exist username"
$ws1.Range("D5").Value = "Both field invalid"
$ws1.Range("D6").Value = "Empty password"
$ws1.Range("D7").Value = "Empty username"
$ws1.Range("D8").Value = "Both field empty"

# ============ Sheet2 (add_user) : header row ============
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "new_password"
$ws2.Range("C1").Value = "fname"
$ws2.Range("D1").Value = "lname"
$ws2.Range("E1").Value = "email"
$ws2.Range("F1").Value = "valid"
$ws2.Range("G1").Value = "testcase name"

# Row 2 - valid testcase
$ws2.Range("A2").Value = "tester"
$ws2.Range("B2").Value = 123456789
$ws2.Range("C2").Value = "Bao"
$ws2.Range("D2").Value = "Tran"
$ws2.Range("E2").Value = "bao@gmail.com"
$ws2.Range("F2").Value = $true
$ws2.Range("G2").Value = "Valid testcase"

# Row 3 - invalid first name (C3 left empty)
$ws2.Range("A3").Value = "tester"
$ws2.Range("B3").Value = 123456789
$ws2.Range("D3").Value = "Tran"
$ws2.Range("E3").Value = "bao@gmail.com"
$ws2.Range("F3").Value = $false
$ws2.Range("G3").Value = "Invalid first name"

# Row 4 - invalid last name (D4 left empty)
$ws2.Range("A4").Value = "tester1"
$ws2.Range("B4").Value = 123456789
$ws2.Range("C4").Value = "Bao"
$ws2.Range("E4").Value = "bao1@gmail.com"
$ws2.Range("F4").Value = $false
$ws2.Range("G4").Value = "Invalid last name"

# Row 5 - invalid email
$ws2.Range("A5").Value = "tester1"
$ws2.Range("B5").Value = 123456789
$ws2.Range("C5").Value = "Bao"
$ws2.Range("D5").Value = "Tran"
$ws2.Range("E5").Value = "admin@school.a"
$ws2.Range("F5").Value = $false
$ws2.Range("G5").Value = "Invalid email"

# Row 6 - invalid password (A6 present, B6 empty)
$ws2.Range("A6").Value = "tester1"
$ws2.Range("C6").Value = "Bao"
$ws2.Range("D6").Value = "Tran"
$ws2.Range("E6").Value = "bao1@gmail.com"
$ws2.Range("F6").Value = $false
$ws2.Range("G6").Value = "Invalid password"

# Row 7 - invalid username
$ws2.Range("A7").Value = "manager"
$ws2.Range("B7").Value = 123456789
$ws2.Range("C7").Value = "Bao"
$ws2.Range("D7").Value = "Tran"
$ws2.Range("E7").Value = "bao@gmail.com"
$ws2.Range("F7").Value = $false
$ws2.Range("G7").Value = "Invalid username"
